$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update D-column "Attribute" cells to add the new "Alle Daten hochgeladen?" upload field
$ws.Range("D4").Value = "HINWEIS --> AMA Daten werden von der Finanz verwertet, saubere Erstdaten helfen bei Plausibilitätsprüfung:info;Flächenaufstellung:select(Eigen-,Pacht,Mitbewirtschaftung):pflicht;Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D5").Value = "Tierarten-und-Bestände:text;Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D6").Value = "Anzahl Hektar:number;typische Nutzung:text:pflicht;Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D8").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D9").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D10").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D11").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D12").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D13").Value = "Beschreibung:text;Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D15").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D16").Value = "Art:select(Neubau,Umbau,Zubau);Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D17").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D19").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D20").Value = "AMA:select(keine,AMA,ÖPUL,Sonstige):pflicht;Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D21").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D22").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D23").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D24").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D29").Value = "Betriebsführer:text;Mitunternehmer:text;Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D30").Value = "Erwerbsart:select(Vollerwerb,Nebenerwerb);Bewirtschaftungsart:select(Bio,konventionell);Kurzbeschreibung:text;Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D31").Value = "Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"
$ws.Range("D32").Value = "Beschreibung:text:pflicht;Alle Daten hochgeladen?:info;Upload:checkbox:pflicht"

# New E16 cell ("Upload" column) gets "Ja"
$ws.Range("E16").Value = "Ja"

# Update sheet view (scroll position / selection) to match final state
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D34").Select()
